$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.309.68"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.031.79"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'578.14"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "'168.68"
$ws.Range("E6").Value = "  +3.50%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.030.01"
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("E12").Value = "  +7.52%  "
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "'36.87"
$ws.Range("E14").Value = "  +6.62%  "
$ws.Range("D16").Value = "66.324.59"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "3.534.00"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("E18").Value = "  +4.63%  "
$ws.Range("D19").Value = "'16.50"
$ws.Range("E19").Value = "  +19.17%  "
$ws.Range("D20").Value = "3.029.16"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "'467.89"
$ws.Range("E21").Value = "  +3.20%  "
$ws.Range("E22").Value = "  +3.56%  "
$ws.Range("D23").Value = "'7.40"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "'83.09"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("D27").Value = "'10.04"
$ws.Range("E27").Value = "  -3.36%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +7.10%  "
$ws.Range("D33").Value = "0.0₃0996"
$ws.Range("E33").Value = "  -4.49%  "
$ws.Range("E34").Value = "  +3.31%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'5.86"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").Value = "'0.990"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'48.19"
$ws.Range("E38").Value = "  +8.66%  "
$ws.Range("D39").Value = "'2.06"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "'49.53"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'8.66"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.85"
$ws.Range("E44").Value = "  -3.78%  "
$ws.Range("D45").Value = "'0.0361"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").Value = "'379.33"
$ws.Range("E46").Value = "  -5.07%  "
$ws.Range("D47").Value = "2.705.00"
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("D48").Value = "'134.27"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'24.52"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("E51").Value = "  +4.34%  "
